# Append the new daily allocation row (09/17/2025) to Sheet1.
# Row 16: Date (as literal text, matching the existing Date column's
# inline-string cells) + BTC / KAS allocation fractions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Date cell to be stored as text rather than letting Excel's
# automatic date-recognition turn "09/17/2025" into a date serial number,
# then clear the temporary number-format override so the cell is left
# with the default (unstyled) formatting, consistent with the rest of
# the Date column.
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "09/17/2025"
$ws.Range("A16").ClearFormats()

$ws.Range("B16").Value = 0.1245762204778458
$ws.Range("C16").Value = 0.8754237795221542
